# Update "想去人数" (F column) values across sheets to reflect refreshed
# crawl output, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 153
$ws1.Range("F4").Value  = 1161
$ws1.Range("F5").Value  = 1060
$ws1.Range("F6").Value  = 1851
$ws1.Range("F7").Value  = 584
$ws1.Range("F8").Value  = 1227
$ws1.Range("F9").Value  = 62
$ws1.Range("F12").Value = 316
$ws1.Range("F13").Value = 93
$ws1.Range("F15").Value = 743
$ws1.Range("F16").Value = 205
$ws1.Range("F21").Value = 179
$ws1.Range("F22").Value = 686
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 656
$ws1.Range("F27").Value = 889
$ws1.Range("F28").Value = 331
$ws1.Range("F29").Value = 173
$ws1.Range("F34").Value = 414

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 326

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 319

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 319
$ws4.Range("F3").Value  = 153
$ws4.Range("F5").Value  = 1161
$ws4.Range("F6").Value  = 1060
$ws4.Range("F7").Value  = 1851
$ws4.Range("F8").Value  = 584
$ws4.Range("F9").Value  = 1227
$ws4.Range("F10").Value = 62
$ws4.Range("F14").Value = 316
$ws4.Range("F15").Value = 93
$ws4.Range("F17").Value = 743
$ws4.Range("F18").Value = 205
$ws4.Range("F22").Value = 326
$ws4.Range("F29").Value = 179
$ws4.Range("F30").Value = 686
$ws4.Range("F31").Value = 56
$ws4.Range("F32").Value = 656
$ws4.Range("F35").Value = 889
$ws4.Range("F36").Value = 331
$ws4.Range("F39").Value = 173
$ws4.Range("F48").Value = 414
